# Update the "K" column (column G) values in Sheet1, rows 2-37.
# These values represent the number of strikeouts (K) and were
# regenerated, replacing the prior "Strike#" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 5
    3  = 3
    4  = 1
    5  = 4
    6  = 1
    7  = 4
    8  = 6
    9  = 3
    10 = 3
    11 = 10
    12 = 3
    13 = 4
    14 = 3
    15 = 4
    16 = 1
    17 = 5
    18 = 5
    19 = 7
    20 = 7
    21 = 2
    22 = 5
    23 = 4
    24 = 5
    25 = 4
    26 = 8
    27 = 3
    28 = 6
    29 = 3
    30 = 6
    31 = 2
    32 = 6
    33 = 4
    34 = 4
    35 = 4
    36 = 5
    37 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
